$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 1612
$ws.Range("F5").Value = 9234
$ws.Range("F6").Value = 281
$ws.Range("F8").Value = 1283
$ws.Range("F10").Value = 669
$ws.Range("F13").Value = 164
$ws.Range("F14").Value = 301
$ws.Range("F17").Value = 1546
$ws.Range("F18").Value = 1337
$ws.Range("F20").Value = 55
$ws.Range("F21").Value = 1414
$ws.Range("F22").Value = 94
$ws.Range("F23").Value = 255
$ws.Range("F25").Value = 105
$ws.Range("F26").Value = 77
$ws.Range("F28").Value = 330
$ws.Range("F29").Value = 330
$ws.Range("F33").Value = 240
$ws.Range("F34").Value = 223
$ws.Range("F35").Value = 63
$ws.Range("F36").Value = 586
$ws.Range("F37").Value = 616
$ws.Range("F38").Value = 430
$ws.Range("F40").Value = 78
$ws.Range("F42").Value = 110
$ws.Range("F43").Value = 529
$ws.Range("F45").Value = 701
$ws.Range("F46").Value = 242
$ws = $wb.Worksheets.Item(2)
$ws.Range("F6").Value = 59
$ws.Range("G8").Value = "不可售"
$ws.Range("F16").Value = 676
$ws.Range("F24").Value = 943
$ws.Range("F26").Value = 1041
$ws.Range("F30").Value = 258
$ws = $wb.Worksheets.Item(3)
$ws.Range("F4").Value = 757
$ws.Range("F5").Value = 337
$ws.Range("F6").Value = 151
$ws.Range("F7").Value = 2231
$ws.Range("F8").Value = 3324
$ws.Range("F9").Value = 39
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 1612
$ws.Range("F4").Value = 757
$ws.Range("F5").Value = 9234
$ws.Range("F6").Value = 337
$ws.Range("F7").Value = 151
$ws.Range("F8").Value = 2231
$ws.Range("F9").Value = 3324
$ws.Range("F11").Value = 1283
$ws.Range("F12").Value = 669
$ws.Range("F15").Value = 164
$ws.Range("F16").Value = 301
$ws.Range("F17").Value = 1546
$ws.Range("F18").Value = 676
$ws.Range("F19").Value = 1337
$ws.Range("F21").Value = 39
$ws.Range("F22").Value = 1414
$ws.Range("F23").Value = 94
$ws.Range("F24").Value = 255
$ws.Range("F26").Value = 105
$ws.Range("F27").Value = 77
$ws.Range("F28").Value = 330
$ws.Range("F29").Value = 330
$ws.Range("F35").Value = 943
$ws.Range("F36").Value = 241
$ws.Range("F38").Value = 223
$ws.Range("F39").Value = 1041
$ws.Range("F41").Value = 586
$ws.Range("F42").Value = 616
$ws.Range("F45").Value = 258
$ws.Range("F47").Value = 111
$ws.Range("F49").Value = 529
$ws.Range("F50").Value = 701
